$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 0.1450048780487805
$ws.Range("V2").Value = 0.0002448603057459146
$ws.Range("Z2").Value = -0.2387391389365841
$ws.Range("AB2").Value = -975.0013919541442
$ws.Range("AC2").Value = "umolO2/min/m2"
$ws.Range("AD2").Value = -975.0013919541442

# Row 3
$ws.Range("T3").Value = 0.1492487804878049
$ws.Range("V3").Value = 0.0001488973818309612
$ws.Range("Z3").Value = -0.2578282185186086
$ws.Range("AB3").Value = -1731.58329144641
$ws.Range("AC3").Value = "umolO2/min/m2"
$ws.Range("AD3").Value = -1731.58329144641

# Row 4
$ws.Range("T4").Value = 0.1469268292682927
$ws.Range("V4").Value = 0.0002222807942365138
$ws.Range("Z4").Value = -0.1827411976190429
$ws.Range("AB4").Value = -822.1187001185556
$ws.Range("AC4").Value = "umolO2/min/m2"
$ws.Range("AD4").Value = -822.1187001185556

# Row 5
$ws.Range("T5").Value = 0.1418926829268293
$ws.Range("V5").Value = 0.0002529432437181515
$ws.Range("Z5").Value = -0.265747140039218
$ws.Range("AB5").Value = -1050.61964151663
$ws.Range("AC5").Value = "umolO2/min/m2"
$ws.Range("AD5").Value = -1050.61964151663

# Row 6
$ws.Range("T6").Value = 0.1446439024390244
$ws.Range("V6").Value = 0.0001851607801792304
$ws.Range("Z6").Value = -0.273584090574137
$ws.Range("AB6").Value = -1477.54881087299
$ws.Range("AC6").Value = "umolO2/min/m2"
$ws.Range("AD6").Value = -1477.54881087299

# Row 7
$ws.Range("T7").Value = 0.1429268292682927
$ws.Range("V7").Value = 0.0003232296608680373
$ws.Range("Z7").Value = -0.2413256967761488
$ws.Range("AB7").Value = -746.6075239755709
$ws.Range("AC7").Value = "umolO2/min/m2"
$ws.Range("AD7").Value = -746.6075239755709

# Row 8
$ws.Range("T8").Value = 0.1544
$ws.Range("V8").Value = 0
$ws.Range("Z8").Value = 0
$ws.Range("AB8").ClearContents()
$ws.Range("AC8").Value = "umolO2/min/m2"
$ws.Range("AD8").ClearContents()

# Row 9
$ws.Range("T9").Value = 0.1450048780487805
$ws.Range("V9").Value = 0.0002448603057459146
$ws.Range("Z9").Value = 0.2455308305614382
$ws.Range("AB9").Value = 1002.738397362859
$ws.Range("AC9").Value = "umolO2/min/m2"
$ws.Range("AD9").Value = 1002.738397362859

# Row 10
$ws.Range("T10").Value = 0.1492487804878049
$ws.Range("V10").Value = 0.0001488973818309612
$ws.Range("Z10").Value = 0.3307141318614385
$ws.Range("AB10").Value = 2221.087622862896
$ws.Range("AC10").Value = "umolO2/min/m2"
$ws.Range("AD10").Value = 2221.087622862896

# Row 11
$ws.Range("T11").Value = 0.1469268292682927
$ws.Range("V11").Value = 0.0002222807942365138
$ws.Range("Z11").Value = 0.1478845783397708
$ws.Range("AB11").Value = 665.3052453214507
$ws.Range("AC11").Value = "umolO2/min/m2"
$ws.Range("AD11").Value = 665.3052453214507

# Row 12
$ws.Range("T12").Value = 0.1418926829268293
$ws.Range("V12").Value = 0.0002529432437181515
$ws.Range("Z12").Value = 0.3020050892897799
$ws.Range("AB12").Value = 1193.963850745493
$ws.Range("AC12").Value = "umolO2/min/m2"
$ws.Range("AD12").Value = 1193.963850745493

# Row 13
$ws.Range("T13").Value = 0.1446439024390244
$ws.Range("V13").Value = 0.0001851607801792304
$ws.Range("Z13").Value = 0.3303973570296264
$ws.Range("AB13").Value = 1784.380886221214
$ws.Range("AC13").Value = "umolO2/min/m2"
$ws.Range("AD13").Value = 1784.380886221214

# Row 14
$ws.Range("T14").Value = 0.1429268292682927
$ws.Range("V14").Value = 0.0003232296608680373
$ws.Range("Z14").Value = 0.2239874616006679
$ws.Range("AB14").Value = 692.9669170804027
$ws.Range("AC14").Value = "umolO2/min/m2"
$ws.Range("AD14").Value = 692.9669170804027

# Row 15
$ws.Range("T15").Value = 0.1544
$ws.Range("V15").Value = 0
$ws.Range("Z15").Value = 0.001216064704046756
$ws.Range("AB15").Value = "Inf"
$ws.Range("AC15").Value = "umolO2/min/m2"
$ws.Range("AD15").Value = "Inf"
